$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 549.5
$ws.Range("I33").Value = 507.66666
$ws.Range("J33").Value = 675
$ws.Range("K33").Value = 507.66666
$ws.Range("L33").Value = 675
$ws.Range("M33").Value = -278.66666
$ws.Range("N33").Value = -1133

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6599.6
$ws.Range("I64").Value = 6749.5
$ws.Range("K64").Value = 6749.5
$ws.Range("M64").Value = -6501.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6599.6
$ws.Range("I67").Value = 6749.5
$ws.Range("K67").Value = 6749.5
$ws.Range("M67").Value = -5891.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 50003490
$ws.Range("I76").Value = 100003120
$ws.Range("J76").Value = 3859.4
$ws.Range("K76").Value = 100003120
$ws.Range("L76").Value = 3859.4
$ws.Range("M76").Value = -100002805
$ws.Range("N76").Value = -4489.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 50003490
$ws.Range("I79").Value = 100003120
$ws.Range("J79").Value = 3859.4
$ws.Range("K79").Value = 100003120
$ws.Range("L79").Value = 3859.4
$ws.Range("M79").Value = -100002028
$ws.Range("N79").Value = -6043.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1177.75
$ws.Range("I98").Value = 863.8333
$ws.Range("J98").Value = 4003
$ws.Range("K98").Value = 863.8333
$ws.Range("L98").Value = 4003
$ws.Range("M98").Value = 634.1667
$ws.Range("N98").Value = -6999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 4833.3335
$ws.Range("I100").Value = 5225.25
$ws.Range("K100").Value = 5225.25
$ws.Range("M100").Value = -4684.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 143.55556
$ws.Range("I107").Value = 143.55556
$ws.Range("K107").Value = 143.55556
$ws.Range("M107").Value = 1776.44444

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1177.75
$ws.Range("I122").Value = 863.8333
$ws.Range("J122").Value = 4003
$ws.Range("K122").Value = 2591.4999
$ws.Range("L122").Value = 12009
$ws.Range("M122").Value = -141.4998999999998
$ws.Range("N122").Value = -16909

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 72433.8
$ws.Range("J133").Value = 72433.8
$ws.Range("L133").Value = 72433.8
$ws.Range("N133").Value = -82553.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 54992.92
$ws.Range("J134").Value = 54992.92
$ws.Range("L134").Value = 54992.92
$ws.Range("N134").Value = -65132.92

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 77983
$ws.Range("J136").Value = 77983
$ws.Range("L136").Value = 77983
$ws.Range("N136").Value = -88183

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1569.9166
$ws.Range("J138").Value = 2494.5
$ws.Range("L138").Value = 7483.5
$ws.Range("N138").Value = -17763.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 98399.78
$ws.Range("J139").Value = 98399.78
$ws.Range("L139").Value = 98399.78
$ws.Range("N139").Value = -108679.78

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 80770.8
$ws.Range("J140").Value = 80770.8
$ws.Range("L140").Value = 80770.8
$ws.Range("N140").Value = -91130.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4028.1743
$ws.Range("I32").Value = 1565.2715
$ws.Range("K32").Value = 1565.2715
$ws.Range("M32").Value = -1278.2715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 51846.85
$ws.Range("I61").Value = 1714.875
$ws.Range("K61").Value = 1714.875
$ws.Range("M61").Value = -1502.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1819.1666
$ws.Range("I63").Value = 1772.1111
$ws.Range("K63").Value = 1772.1111
$ws.Range("M63").Value = -1086.1111

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1819.1666
$ws.Range("I66").Value = 1772.1111
$ws.Range("K66").Value = 8860.5555
$ws.Range("M66").Value = -5428.5555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 51846.85
$ws.Range("I136").Value = 1714.875
$ws.Range("K136").Value = 5144.625
$ws.Range("M136").Value = -2594.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 81241.234
$ws.Range("I22").Value = 104848.7
$ws.Range("J22").Value = 2549.6667
$ws.Range("K22").Value = 104848.7
$ws.Range("L22").Value = 2549.6667
$ws.Range("M22").Value = -104675.7
$ws.Range("N22").Value = -2895.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 3124.25
$ws.Range("J80").Value = 3124.25
$ws.Range("L80").Value = 3124.25
$ws.Range("N80").Value = -5120.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 3124.25
$ws.Range("J83").Value = 3124.25
$ws.Range("L83").Value = 15621.25
$ws.Range("N83").Value = -25605.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 146714
$ws.Range("I105").Value = 335999.34
$ws.Range("K105").Value = 335999.34
$ws.Range("M105").Value = -334252.34

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 34821.777
$ws.Range("J132").Value = 34821.777
$ws.Range("L132").Value = 34821.777
$ws.Range("N132").Value = -44941.777

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 105163.336
$ws.Range("J135").Value = 105163.336
$ws.Range("L135").Value = 105163.336
$ws.Range("N135").Value = -115303.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 99754.86
$ws.Range("J138").Value = 99754.86
$ws.Range("L138").Value = 99754.86
$ws.Range("N138").Value = -110034.86

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2217.6572
$ws.Range("J31").Value = 4311.7144
$ws.Range("L31").Value = 4311.7144
$ws.Range("N31").Value = -4901.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2217.6572
$ws.Range("J34").Value = 4311.7144
$ws.Range("L34").Value = 4311.7144
$ws.Range("N34").Value = -4715.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2068953.6
$ws.Range("I132").Value = 2843687
$ws.Range("K132").Value = 8531061
$ws.Range("M132").Value = -8528531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 94496
$ws.Range("J138").Value = 94496
$ws.Range("L138").Value = 94496
$ws.Range("N138").Value = -104776

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 875.5
$ws.Range("I107").Value = 895.2
$ws.Range("J107").Value = 777
$ws.Range("K107").Value = 2685.6
$ws.Range("L107").Value = 2331
$ws.Range("M107").Value = -765.6000000000004
$ws.Range("N107").Value = -6171

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 62249.188
$ws.Range("J70").Value = 154080.67
$ws.Range("L70").Value = 154080.67
$ws.Range("N70").Value = -154620.67

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 62249.188
$ws.Range("J73").Value = 154080.67
$ws.Range("L73").Value = 154080.67
$ws.Range("N73").Value = -155952.67

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 4500
$ws.Range("K80").Value = 3500
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -2502
$ws.Range("N80").Value = -6496

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 4500
$ws.Range("K83").Value = 17500
$ws.Range("L83").Value = 22500
$ws.Range("M83").Value = -12508
$ws.Range("N83").Value = -32484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 30266.363
$ws.Range("J109").Value = 30266.363
$ws.Range("L109").Value = 30266.363
$ws.Range("N109").Value = -32346.363

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2024.3334
$ws.Range("J113").Value = 2798.6
$ws.Range("L113").Value = 2798.6
$ws.Range("N113").Value = -7138.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5117
$ws.Range("I132").Value = 4280.222
$ws.Range("K132").Value = 12840.666
$ws.Range("M132").Value = -10310.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 69996.664
$ws.Range("J135").Value = 69996.664
$ws.Range("L135").Value = 69996.664
$ws.Range("N135").Value = -80136.664

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 90068.57000000001
$ws.Range("J140").Value = 89996.664
$ws.Range("L140").Value = 89996.664
$ws.Range("N140").Value = -100356.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1434.3235
$ws.Range("I136").Value = 1235.1482
$ws.Range("K136").Value = 3705.4446
$ws.Range("M136").Value = -1155.4446
